$wb = $excel.ActiveWorkbook

# --- Text change: "Ready for handoff" -> "In Translation" ---
# Overview sheet: status columns are zh-cn (E) and de-de (F) for rows 2-4
$wsOverview = $wb.Sheets.Item("Overview")
$wsOverview.Range("E2:F4").Value = "In Translation"

# zh-cn sheet: Status column is C, rows 2-4
$wsZhCn = $wb.Sheets.Item("zh-cn")
$wsZhCn.Range("C2:C4").Value = "In Translation"

# de-de sheet: Status column is C, rows 2-4
$wsDeDe = $wb.Sheets.Item("de-de")
$wsDeDe.Range("C2:C4").Value = "In Translation"

# --- Column width changes (narrower status columns) ---
# Overview: columns E and F (zh-cn / de-de status) get narrower
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# zh-cn: column C (Status) gets narrower
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

# de-de: column C (Status) gets narrower
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
